# Mark the remaining "Backend Models" design tasks as Done, widen the
# Task Type description column, and move the viewport/selection down to
# where the newly completed rows live.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status column (F): these rows moved from "Not Started" to "Done"
$doneRows = @(15, 17, 22, 25, 27, 28, 29, 30, 31)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 6).Value = "Done"
}

# Column D ("Task Type" details) got wider to fit the longer class/method
# descriptions.
$ws.Columns.Item(4).ColumnWidth = 62.65

# Scroll the view down and select F31, matching where work left off.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F31").Select()
